# Publication preparation for release 0.2.0
# - bump the Version metadata value
# - bump the Date metadata value
# - insert a new "Jurisdiction" row (with its ISO country-code value)
#   right after the "Contact" row, pushing the remaining metadata rows
#   (Description ... Count) down by one

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Shift rows 11..21 down to 12..22, one row at a time, bottom-up, so the
# ranges never overlap. A plain cell Cut/Paste (rather than re-typing the
# .Value) keeps each cell's original type/format exactly as stored.
for ($r = 21; $r -ge 11; $r--) {
    $ws.Range("A" + $r + ":B" + $r).Cut($ws.Range("A" + ($r + 1) + ":B" + ($r + 1)))
}

# New "Jurisdiction" row now occupies row 11 (it inherited row 11's prior
# formatting from the cut/shift above, matching the surrounding rows)
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"

# Bump Version and Date values for the 0.2.0 publication
$ws.Range("B3").Value = "0.2.0"
$ws.Range("B8").Value = "2023-10-19T17:05:12+00:00"
